# no-op test
$p = $ppt.ActivePresentation
Write-Host "slides:" $p.Slides.Count
